$d = $word.ActiveDocument

# 1. Rename the document title heading.
$d.Content.Find.Execute(
    "User Stories the application demonstrates", $true, $false, $false, $false, $false,
    $true, 1, $false, "User Stories Implemented", 2)

# 2. Drop the "Notes" section (its heading plus the explanatory paragraph
#    that follows it) that used to sit right before the "How to run the
#    application" heading. Deleting the range that spans from the start of
#    the "Notes" paragraph's text up to (but not including) the "How to run
#    the application" text merges what is left of the "Notes" paragraph
#    straight into the following heading paragraph, exactly as the diff
#    shows.
$full = $d.Content.Text
$startIdx = $full.IndexOf("Notes ")
$endIdx = $full.IndexOf("How to run the application")
if ($startIdx -ge 0 -and $endIdx -gt $startIdx) {
    $d.Range($startIdx, $endIdx).Delete()
}
